$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency Price (D) and Volume/1h (E) columns with the latest scraped
# values for the Jan 29, 2023 data run. Values are stored as literal text, so the
# cell is temporarily switched to a Text number format while the value is assigned,
# then its original style is restored so formatting is left unchanged.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "309.02"
Set-TextValue $ws.Range("E2") "-0.51%"
Set-TextValue $ws.Range("D3") "39.76"
Set-TextValue $ws.Range("E3") "2.25%"
Set-TextValue $ws.Range("D4") "5.130"
Set-TextValue $ws.Range("E4") "0.14%"
Set-TextValue $ws.Range("D5") "0.08140"
Set-TextValue $ws.Range("D6") "1.948"
Set-TextValue $ws.Range("E6") "-2.71%"
Set-TextValue $ws.Range("D7") "8.118"
Set-TextValue $ws.Range("E7") "2.45%"
Set-TextValue $ws.Range("D8") "0.9286"
Set-TextValue $ws.Range("E8") "-0.42%"
Set-TextValue $ws.Range("D9") "0.1418"
Set-TextValue $ws.Range("E9") "0.85%"
Set-TextValue $ws.Range("D10") "0.1931"
Set-TextValue $ws.Range("E10") "-1.30%"
Set-TextValue $ws.Range("D11") "0.09077"
Set-TextValue $ws.Range("E11") "-0.69%"
Set-TextValue $ws.Range("D12") "0.03509"
Set-TextValue $ws.Range("E12") "1.15%"
Set-TextValue $ws.Range("D13") "0.09806"
Set-TextValue $ws.Range("E13") "-0.45%"
Set-TextValue $ws.Range("D14") "0.001394"
Set-TextValue $ws.Range("E14") "-1.22%"
Set-TextValue $ws.Range("D15") "0.005859"
Set-TextValue $ws.Range("E15") "0.58%"
Set-TextValue $ws.Range("D16") "3.917"
Set-TextValue $ws.Range("E16") "9.78%"
Set-TextValue $ws.Range("D17") "4.218"
Set-TextValue $ws.Range("E17") "0.50%"
Set-TextValue $ws.Range("D19") "0.3454"
Set-TextValue $ws.Range("E19") "0.12%"
Set-TextValue $ws.Range("D20") "0.1311"
Set-TextValue $ws.Range("E20") "-2.03%"
Set-TextValue $ws.Range("D21") "4.722"
Set-TextValue $ws.Range("E21") "-2.06%"
Set-TextValue $ws.Range("D22") "0.2424"
Set-TextValue $ws.Range("E22") "-1.84%"
Set-TextValue $ws.Range("D23") "0.04382"
Set-TextValue $ws.Range("E23") "-1.98%"
Set-TextValue $ws.Range("D24") "0.001232"
Set-TextValue $ws.Range("E24") "-0.59%"
Set-TextValue $ws.Range("D25") "0.004385"
Set-TextValue $ws.Range("E25") "5.08%"
Set-TextValue $ws.Range("D26") "0.0001301"
Set-TextValue $ws.Range("E26") "-0.11%"
Set-TextValue $ws.Range("E27") "-9.97%"
Set-TextValue $ws.Range("D39") "0.02064"
Set-TextValue $ws.Range("E39") "-2.41%"
Set-TextValue $ws.Range("D40") "0.05103"
Set-TextValue $ws.Range("E40") "-1.59%"
Set-TextValue $ws.Range("D41") "0.007432"
Set-TextValue $ws.Range("E41") "-0.63%"
Set-TextValue $ws.Range("D42") "0.009871"
Set-TextValue $ws.Range("E42") "-1.41%"
Set-TextValue $ws.Range("E43") "-0.20%"
Set-TextValue $ws.Range("D44") "0.002132"
Set-TextValue $ws.Range("E44") "-0.11%"
Set-TextValue $ws.Range("D45") "0.009566"
Set-TextValue $ws.Range("E45") "-2.09%"
Set-TextValue $ws.Range("D46") "0.00006382"
Set-TextValue $ws.Range("E46") "0.87%"
Set-TextValue $ws.Range("D47") "0.00000000751"
Set-TextValue $ws.Range("E47") "-0.12%"
Set-TextValue $ws.Range("D48") "0.002716"
Set-TextValue $ws.Range("E49") "-18.90%"
Set-TextValue $ws.Range("D50") "0.00002102"
Set-TextValue $ws.Range("E50") "-0.12%"
Set-TextValue $ws.Range("D51") "0.0002002"
Set-TextValue $ws.Range("E51") "-0.12%"
